$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Form the consolidated "Absent" column (H) values
$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 0
